$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two now-unused rows at the bottom (data shifted up; min-count
# filtering now drops what used to be rows 36-37)
$ws.Rows("36:37").Delete()

$ws.Range("B3").Value = 0.9782608695652174
$ws.Range("C3").Value = 45
$ws.Range("D3").Value = 45
$ws.Range("H3").Value = 1
$ws.Range("K3").Value = 0.8392857142857143
$ws.Range("L3").Value = 47
$ws.Range("M3").Value = 47
$ws.Range("Q3").Value = 9
$ws.Range("B4").Value = 0.8863636363636364
$ws.Range("C4").Value = 39
$ws.Range("D4").Value = 39
$ws.Range("H4").Value = 5
$ws.Range("K4").Value = 0.8307692307692308
$ws.Range("L4").Value = 54
$ws.Range("M4").Value = 54
$ws.Range("Q4").Value = 11
$ws.Range("A5").Value = "poor"
$ws.Range("B5").Value = 0.7887323943661971
$ws.Range("C5").Value = 56
$ws.Range("D5").Value = 56
$ws.Range("K5").Value = 0.6774193548387096
$ws.Range("L5").Value = 63
$ws.Range("M5").Value = 63
$ws.Range("Q5").Value = 30
$ws.Range("A6").Value = "however"
$ws.Range("B6").Value = 0.765625
$ws.Range("C6").Value = 49
$ws.Range("D6").Value = 49
$ws.Range("H6").Value = 15
$ws.Range("B7").Value = 0.7475728155339806
$ws.Range("C7").Value = 154
$ws.Range("D7").Value = 154
$ws.Range("H7").Value = 52
$ws.Range("K7").Value = 0.53125
$ws.Range("L7").Value = 34
$ws.Range("M7").Value = 34
$ws.Range("Q7").Value = 30
$ws.Range("A8").Value = "disappointed"
$ws.Range("B8").Value = 0.7473118279569892
$ws.Range("C8").Value = 139
$ws.Range("D8").Value = 139
$ws.Range("H8").Value = 47
$ws.Range("A9").Value = "waste"
$ws.Range("B9").Value = 0.6283783783783784
$ws.Range("C9").Value = 93
$ws.Range("D9").Value = 93
$ws.Range("H9").Value = 55
$ws.Range("K9").Value = 0.3795081967213115
$ws.Range("L9").Value = 463
$ws.Range("M9").Value = 463
$ws.Range("Q9").Value = 757
$ws.Range("A10").Value = "guess"
$ws.Range("B10").Value = 0.5925925925925926
$ws.Range("C10").Value = 32
$ws.Range("D10").Value = 32
$ws.Range("H10").Value = 22
$ws.Range("K10").Value = 0.3486370157819225
$ws.Range("L10").Value = 243
$ws.Range("M10").Value = 243
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 454
$ws.Range("A11").Value = "smaller"
$ws.Range("B11").Value = 0.5798319327731093
$ws.Range("C11").Value = 69
$ws.Range("D11").Value = 69
$ws.Range("H11").Value = 50
$ws.Range("K11").Value = 0.3112033195020747
$ws.Range("L11").Value = 150
$ws.Range("M11").Value = 150
$ws.Range("Q11").Value = 332
$ws.Range("A12").Value = "junk"
$ws.Range("B12").Value = 0.5636363636363636
$ws.Range("C12").Value = 31
$ws.Range("D12").Value = 31
$ws.Range("H12").Value = 24
$ws.Range("J12").Value = "best"
$ws.Range("K12").Value = 0.2583333333333334
$ws.Range("L12").Value = 31
$ws.Range("M12").Value = 31
$ws.Range("Q12").Value = 89
$ws.Range("A13").Value = "small"
$ws.Range("B13").Value = 0.518840579710145
$ws.Range("C13").Value = 179
$ws.Range("D13").Value = 179
$ws.Range("H13").Value = 166
$ws.Range("J13").Value = "perfect"
$ws.Range("K13").Value = 0.2409638554216867
$ws.Range("L13").Value = 40
$ws.Range("M13").Value = 40
$ws.Range("Q13").Value = 126
$ws.Range("A14").Value = "broken"
$ws.Range("B14").Value = 0.4819277108433735
$ws.Range("C14").Value = 40
$ws.Range("D14").Value = 40
$ws.Range("H14").Value = 43
$ws.Range("K14").Value = 0.2063492063492063
$ws.Range("L14").Value = 39
$ws.Range("M14").Value = 39
$ws.Range("Q14").Value = 150
$ws.Range("A15").Value = "plastic"
$ws.Range("B15").Value = 0.4645669291338583
$ws.Range("C15").Value = 59
$ws.Range("D15").Value = 59
$ws.Range("H15").Value = 68
$ws.Range("J15").Value = "loved"
$ws.Range("K15").Value = 0.1865443425076453
$ws.Range("L15").Value = 61
$ws.Range("M15").Value = 61
$ws.Range("Q15").Value = 266
$ws.Range("A16").Value = "apart"
$ws.Range("B16").Value = 0.4526315789473684
$ws.Range("C16").Value = 43
$ws.Range("D16").Value = 43
$ws.Range("H16").Value = 52
$ws.Range("J16").Value = "christmas"
$ws.Range("K16").Value = 0.1244979919678715
$ws.Range("L16").Value = 31
$ws.Range("M16").Value = 31
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = $false
$ws.Range("Q16").Value = 218
$ws.Range("B17").Value = 0.3820224719101123
$ws.Range("C17").Value = 34
$ws.Range("D17").Value = 34
$ws.Range("H17").Value = 55
$ws.Range("J17").Value = "fun"
$ws.Range("K17").Value = 0.1087719298245614
$ws.Range("L17").Value = 124
$ws.Range("M17").Value = 125
$ws.Range("N17").Value = 0.99
$ws.Range("O17").Value = 0.01000000000000001
$ws.Range("P17").Value = $true
$ws.Range("Q17").Value = 1016
$ws.Range("A18").Value = "ok"
$ws.Range("B18").Value = 0.3671875
$ws.Range("C18").Value = 47
$ws.Range("D18").Value = 47
$ws.Range("H18").Value = 81
$ws.Range("J18").Value = "game"
$ws.Range("K18").Value = 0.04935064935064935
$ws.Range("L18").Value = 76
$ws.Range("M18").Value = 77
$ws.Range("N18").Value = 0.99
$ws.Range("O18").Value = 0.01000000000000001
$ws.Range("P18").Value = $true
$ws.Range("Q18").Value = 1464
$ws.Range("A19").Value = "thought"
$ws.Range("B19").Value = 0.3415841584158416
$ws.Range("C19").Value = 69
$ws.Range("D19").Value = 69
$ws.Range("H19").Value = 133
$ws.Range("A20").Value = "cheap"
$ws.Range("B20").Value = 0.3364928909952606
$ws.Range("C20").Value = 71
$ws.Range("D20").Value = 71
$ws.Range("H20").Value = 140
$ws.Range("A21").Value = "though"
$ws.Range("B21").Value = 0.2905982905982906
$ws.Range("C21").Value = 34
$ws.Range("D21").Value = 34
$ws.Range("H21").Value = 83
$ws.Range("A22").Value = "item"
$ws.Range("B22").Value = 0.2318840579710145
$ws.Range("C22").Value = 64
$ws.Range("D22").Value = 64
$ws.Range("H22").Value = 212
$ws.Range("B23").Value = 0.2164948453608248
$ws.Range("C23").Value = 42
$ws.Range("D23").Value = 42
$ws.Range("H23").Value = 152
$ws.Range("A24").Value = "hard"
$ws.Range("B24").Value = 0.215
$ws.Range("C24").Value = 43
$ws.Range("D24").Value = 43
$ws.Range("H24").Value = 157
$ws.Range("B25").Value = 0.2120253164556962
$ws.Range("C25").Value = 67
$ws.Range("D25").Value = 67
$ws.Range("H25").Value = 249
$ws.Range("A26").Value = "work"
$ws.Range("B26").Value = 0.2025316455696203
$ws.Range("C26").Value = 64
$ws.Range("D26").Value = 64
$ws.Range("H26").Value = 252
$ws.Range("A27").Value = "would"
$ws.Range("B27").Value = 0.1810089020771513
$ws.Range("C27").Value = 122
$ws.Range("D27").Value = 122
$ws.Range("H27").Value = 552
$ws.Range("A28").Value = "product"
$ws.Range("B28").Value = 0.1585903083700441
$ws.Range("C28").Value = 72
$ws.Range("D28").Value = 72
$ws.Range("H28").Value = 382
$ws.Range("A29").Value = "better"
$ws.Range("B29").Value = 0.1448598130841121
$ws.Range("C29").Value = 31
$ws.Range("D29").Value = 31
$ws.Range("H29").Value = 183
$ws.Range("B30").Value = 0.1408045977011494
$ws.Range("C30").Value = 49
$ws.Range("D30").Value = 49
$ws.Range("H30").Value = 299
$ws.Range("A31").Value = "3"
$ws.Range("B31").Value = 0.1376518218623482
$ws.Range("C31").Value = 34
$ws.Range("D31").Value = 35
$ws.Range("E31").Value = 0.03
$ws.Range("F31").Value = 0.97
$ws.Range("G31").Value = $true
$ws.Range("H31").Value = 213
$ws.Range("A32").Value = "use"
$ws.Range("B32").Value = 0.09315068493150686
$ws.Range("C32").Value = 34
$ws.Range("D32").Value = 34
$ws.Range("H32").Value = 331
$ws.Range("A33").Value = "little"
$ws.Range("B33").Value = 0.08685968819599109
$ws.Range("C33").Value = 39
$ws.Range("D33").Value = 39
$ws.Range("H33").Value = 410
$ws.Range("A34").Value = "like"
$ws.Range("B34").Value = 0.06919275123558484
$ws.Range("C34").Value = 42
$ws.Range("D34").Value = 43
$ws.Range("E34").Value = 0.02
$ws.Range("F34").Value = 0.98
$ws.Range("G34").Value = $true
$ws.Range("H34").Value = 565
$ws.Range("A35").Value = "one"
$ws.Range("B35").Value = 0.04689480354879594
$ws.Range("C35").Value = 37
$ws.Range("D35").Value = 42
$ws.Range("E35").Value = 0.12
$ws.Range("F35").Value = 0.88
$ws.Range("H35").Value = 752
